$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for cell content/formatting so we don't drag along
# the old "Яндекс лицей / GitHub" course data or its per-cell wrap formatting.
$ws.Cells.Clear()

# --- Row 1: Yandex Lyceum course -------------------------------------------------
$ws.Range("A1").Value = "Курсы"
$ws.Range("B1").Value = "Яндекс Лицей"
$ws.Range("C1").Value = 43344
$ws.Range("C1").NumberFormat = "mmm-yy"
$ws.Range("D1").Value = "Обучение програмированию на языке Python на базе компании Яндекс."
$ws.Range("E1").Value = "icon/yandex.jpg"

# --- Row 2: Super-English course -------------------------------------------------
$ws.Range("A2").Value = "Онлайн обучение"
$ws.Range("B2").Value = "Super-English"
$ws.Range("C2").Value = 43862
$ws.Range("C2").NumberFormat = "mmm-yy"
$ws.Range("D2").Value = "Изучение английского языка с нуля, до свободного общения вместе с Петровой Оксаной Сергеевной."
$ws.Range("E2").Value = "icon/English.jpg"

# Rows should size themselves automatically again (no forced 150/315 height).
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

# Column widths for the new 5-column layout.
$ws.Columns.Item(1).ColumnWidth = 22.833333333333332
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
$ws.Columns.Item(4).ColumnWidth = 26.666666666666668
$ws.Columns.Item(5).ColumnWidth = 17.333333333333332

# Move the active selection like in the saved workbook.
$ws.Range("C5").Select()
